$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that need a "Processed" marker added in column C.
$rows = @(21, 22, 23, 24, 25)
for ($r = 45; $r -le 1035; $r += 10) {
    $rows += $r
}

foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = "Processed"
}
